$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 130.6
$ws.Range("B2").Value = 4.176122603564219
$ws.Range("C2").Value = 135.5
$ws.Range("D2").Value = 3.640054944640259
$ws.Range("E2").Value = 132.7
$ws.Range("F2").Value = 3.950949253027682
$ws.Range("G2").Value = 131.3
$ws.Range("H2").Value = 3.226453160980336
$ws.Range("I2").Value = 131.2
$ws.Range("J2").Value = 4.019950248448356
$ws.Range("K2").Value = 134.3
$ws.Range("L2").Value = 5.273518749374084
$ws.Range("M2").Value = 131
$ws.Range("N2").Value = 3.898717737923585
$ws.Range("O2").Value = 132.5
$ws.Range("P2").Value = 3.263433774416144
$ws.Range("Q2").Value = 128.5
$ws.Range("R2").Value = 4.031128874149275
$ws.Range("S2").Value = 129.9
$ws.Range("T2").Value = 2.3
$ws.Range("U2").Value = 125.5
$ws.Range("V2").Value = 3.556683848755748
$ws.Range("W2").Value = 123.9
$ws.Range("X2").Value = 4.784349485562274
$ws.Range("Y2").Value = 124.1
$ws.Range("Z2").Value = 3.176476034853718
$ws.Range("AA2").Value = 124.8
$ws.Range("AB2").Value = 3.572114219898351
$ws.Range("AC2").Value = 123.5
$ws.Range("AD2").Value = 4.129164564412516
$ws.Range("AE2").Value = 127.4
$ws.Range("AF2").Value = 3.006659275674582
$ws.Range("AG2").Value = 123.8
$ws.Range("AH2").Value = 5.134199061197374
$ws.Range("AI2").Value = 122.7
$ws.Range("AJ2").Value = 3.551056180912941
$ws.Range("AK2").Value = 121.9
$ws.Range("AL2").Value = 3.448187929913333
$ws.Range("AM2").Value = 121.2
$ws.Range("AN2").Value = 4.833218389437829
$ws.Range("AO2").Value = 123.6
$ws.Range("AP2").Value = 5.730619512757761
$ws.Range("AQ2").Value = 122.6
$ws.Range("AR2").Value = 6.755738301621815
$ws.Range("AS2").Value = 124.5
$ws.Range("AT2").Value = 5.142956348249516
$ws.Range("AU2").Value = 123.4
$ws.Range("AV2").Value = 5.834380858325929
$ws.Range("AW2").Value = 121.9
$ws.Range("AX2").Value = 7.147726911403373
$ws.Range("AY2").Value = 122.3
$ws.Range("AZ2").Value = 8.626123115281858
$ws.Range("BA2").Value = 117.1
$ws.Range("BB2").Value = 3.83275357934736
$ws.Range("BC2").Value = 119.7
$ws.Range("BD2").Value = 7.071774883294858
$ws.Range("BE2").Value = 121.3
$ws.Range("BF2").Value = 9.121951545584968
$ws.Range("BG2").Value = 121.6
$ws.Range("BH2").Value = 10.32666451474047
$ws.Range("BI2").Value = 3774.3
